# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (interest count) and "最低票价" (lowest price)
# figures to the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 525
$ws.Range("F4").Value = 1520
$ws.Range("F8").Value = 152
$ws.Range("F9").Value = 738
$ws.Range("F11").Value = 63
$ws.Range("F13").Value = 52
$ws.Range("F14").Value = 6398
$ws.Range("F15").Value = 8
$ws.Range("F20").Value = 15308
$ws.Range("G20").Value = 19.9
$ws.Range("F21").Value = 1521
$ws.Range("F23").Value = 142
$ws.Range("F25").Value = 11046
$ws.Range("F27").Value = 4319
$ws.Range("F28").Value = 236

# --- Sheet: 全部类型 ---
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F3").Value = 525
$ws2.Range("F4").Value = 1520
$ws2.Range("F9").Value = 152
$ws2.Range("F10").Value = 738
$ws2.Range("F13").Value = 63
$ws2.Range("F15").Value = 52
$ws2.Range("F17").Value = 6398
$ws2.Range("F18").Value = 8
$ws2.Range("F23").Value = 15308
$ws2.Range("G23").Value = 19.9
$ws2.Range("F24").Value = 1521
$ws2.Range("F26").Value = 142
$ws2.Range("F28").Value = 11046
$ws2.Range("F30").Value = 4319
$ws2.Range("F31").Value = 236

$wb.Save()
